# Apply the crypto price/volume refresh captured in the commit:
#   "Updated cryptos list on Fri May 12 17:35:12 UTC 2023 with GitHub Actions"
#
# Columns B (Coin) / C (Link) / D (Price) / E (Volume(1h)) on the single sheet.
# D holds numeric-looking text (e.g. "26.372.24", "0.3596") that must stay TEXT,
# not be reinterpreted as numbers/dates. A leading apostrophe forces Excel's
# quote-prefix (text) entry mode for those; it is stripped from the stored value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.372.24'
$ws.Range("E2").Value = '  -2.68%  '

# Row 3
$ws.Range("D3").Value = '1.775.86'
$ws.Range("E3").Value = '  -1.10%  '

# Row 4
$ws.Range("E4").Value = '  -0.69%  '

# Row 5
$ws.Range("B5").Value = 'USDC'
$ws.Range("C5").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D5").Value = '''1.000'
$ws.Range("E5").Value = '  -0.73%  '

# Row 6
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '''304.75'
$ws.Range("E6").Value = '  -1.28%  '

# Row 7
$ws.Range("E7").Value = '  +1.20%  '

# Row 8
$ws.Range("D8").Value = '''0.3596'
$ws.Range("E8").Value = '  +1.00%  '

# Row 9
$ws.Range("D9").Value = '''0.07145'
$ws.Range("E9").Value = '  +1.03%  '

# Row 10
$ws.Range("D10").Value = '''0.8362'
$ws.Range("E10").Value = '  -0.87%  '

# Row 11
$ws.Range("D11").Value = '''20.40'
$ws.Range("E11").Value = '  +1.70%  '

# Row 12
$ws.Range("D12").Value = '1.763.91'
$ws.Range("E12").Value = '  -3.07%  '

# Row 13
$ws.Range("D13").Value = '''6.452'
$ws.Range("E13").Value = '  +1.84%  '

# Row 14
$ws.Range("E14").Value = '  -0.68%  '

# Row 15
$ws.Range("D15").Value = '''0.06861'
$ws.Range("E15").Value = '  +0.14%  '

# Row 16
$ws.Range("D16").Value = '''1.005'
$ws.Range("E16").Value = '  -0.42%  '

# Row 17
$ws.Range("E17").Value = '  -0.91%  '

# Row 18
$ws.Range("D18").Value = '''0.000008647'
$ws.Range("E18").Value = '  -0.87%  '

# Row 19
$ws.Range("D19").Value = '''1.000'
$ws.Range("E19").Value = '  -0.72%  '

# Row 20
$ws.Range("D20").Value = '''14.89'
$ws.Range("E20").Value = '  -1.05%  '

# Row 21
$ws.Range("D21").Value = '26.384.98'
$ws.Range("E21").Value = '  -3.17%  '

# Row 22
$ws.Range("D22").Value = '''5.077'
$ws.Range("E22").Value = '  +0.54%  '

# Row 23
$ws.Range("D23").Value = '''10.95'
$ws.Range("E23").Value = '  +2.00%  '

# Row 24
$ws.Range("D24").Value = '1.999.78'
$ws.Range("E24").Value = '  -2.51%  '

# Row 25
$ws.Range("D25").Value = '''151.94'
$ws.Range("E25").Value = '  -0.87%  '

# Row 26
$ws.Range("D26").Value = '''1.801'
$ws.Range("E26").Value = '  -8.11%  '

# Row 27
$ws.Range("E27").Value = '  -0.92%  '

# Row 28
$ws.Range("D28").Value = '''5.068'
$ws.Range("E28").Value = '  +0.90%  '

# Row 29
$ws.Range("D29").Value = '''114.64'
$ws.Range("E29").Value = '  +1.80%  '

# Row 30
$ws.Range("D30").Value = '''1.836'
$ws.Range("E30").Value = '  +10.56%  '

# Row 31
$ws.Range("D31").Value = '''0.08827'
$ws.Range("E31").Value = '  -0.69%  '

# Row 32
$ws.Range("D32").Value = '''0.7262'
$ws.Range("E32").Value = '  +0.22%  '

# Row 33
$ws.Range("D33").Value = '''1.119'
$ws.Range("E33").Value = '  +4.22%  '

# Row 34
$ws.Range("D34").Value = '''4.319'
$ws.Range("E34").Value = '  -1.00%  '

# Row 35
$ws.Range("D35").Value = '''0.9995'
$ws.Range("E35").Value = '  -0.69%  '

# Row 36
$ws.Range("D36").Value = '''2.730'
$ws.Range("E36").Value = '  -5.11%  '

# Row 37
$ws.Range("D37").Value = '''1.090'
$ws.Range("E37").Value = '  +1.83%  '

# Row 38
$ws.Range("E38").Value = '  +0.47%  '

# Row 39
$ws.Range("D39").Value = '''0.01881'
$ws.Range("E39").Value = '  -0.56%  '

# Row 40
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '''0.4914'
$ws.Range("E40").Value = '  -0.57%  '

# Row 41
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '''0.1606'
$ws.Range("E41").Value = '  -0.77%  '

# Row 42
$ws.Range("D42").Value = '''2.607'
$ws.Range("E42").Value = '  -2.18%  '

# Row 43
$ws.Range("D43").Value = '''6.327'
$ws.Range("E43").Value = '  +1.37%  '

# Row 44
$ws.Range("D44").Value = '''7.966'
$ws.Range("E44").Value = '  -0.73%  '

# Row 45
$ws.Range("D45").Value = '''104.53'
$ws.Range("E45").Value = '  -0.11%  '

# Row 46
$ws.Range("D46").Value = '''10.16'
$ws.Range("E46").Value = '  +0.06%  '

# Row 47
$ws.Range("D47").Value = '''0.9993'
$ws.Range("E47").Value = '  -0.76%  '

# Row 48
$ws.Range("E48").Value = '  +2.73%  '

# Row 49
$ws.Range("D49").Value = '''0.06174'
$ws.Range("E49").Value = '  -2.20%  '

# Row 50
$ws.Range("D50").Value = '''0.4445'
$ws.Range("E50").Value = '  -1.92%  '

# Row 51
$ws.Range("E51").Value = '  +3.43%  '
